$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting existing rows 13:114 down to 14:115
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new price-report entry
$ws.Cells.Item(13, 1).Value = 5
$ws.Cells.Item(13, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(13, 3).Value = "Maule"
$ws.Cells.Item(13, 4).Value = 44532
$ws.Cells.Item(13, 5).Value = 7
$ws.Cells.Item(13, 6).Value = 100112031
$ws.Cells.Item(13, 7).Value = "Poroto verde"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 150
$ws.Cells.Item(13, 11).Value = 20000
$ws.Cells.Item(13, 12).Value = 20000
$ws.Cells.Item(13, 13).Value = 20000
$ws.Cells.Item(13, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(13, 15).Value = "Región del Maule"
$ws.Cells.Item(13, 16).Value = 800
$ws.Cells.Item(13, 17).Value = 25
$ws.Cells.Item(13, 18).Value = "Hortaliza"
